# "Generate Report for Handoff"
#
# Refreshes the localization-status report: the source markdown file was
# renamed/regenerated (new GUID-based name + new content hash), a new
# handoff round was generated, and the previous handback (target) files
# were cleared out pending the new handback.

$wb = $excel.ActiveWorkbook

$oldId   = "de2f0e6c-f087-4fa4-9c6f-c804fd303f52"
$newId   = "34b04bb6-7bcc-411c-a239-1c1ec4894557"
$oldHash = "1fc4f996a9086de4e1176e20ec32355e52f818d3"
$newHash = "e8e9e5dc77a745d211f135e5779ec13c124c2016"

$oldFile = "$oldId.md"
$newFile = "$newId.md"
$oldPath = "e2e\$oldId.md"
$newPath = "e2e\$newId.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = $newPath
$ws.Range("G2").Value = "2016-08-18 00:55:59"

$ws.Columns.Item(1).ColumnWidth = 40

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-18 00:55:54"
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(1).ColumnWidth = 40
$ws.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-18 00:55:59"
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(1).ColumnWidth = 40
$ws.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426
